# Update cryptos list data (generated to match GitHub Actions run on 2024-08-17)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 <-> Row 18 swap: WrappedEther moves to row 18, ShibaInu moves to row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'2.575.48"
$ws.Range("E18").Value = "  +1.07%  "

# Price (D) and Volume(1h) (E) updates per row
$ws.Range("D2").Value = "58.928.13"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "2.584.65"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'527.68"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").Value = "'138.40"
$ws.Range("E6").Value = "  -3.47%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("D9").Value = "2.597.17"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").Value = "'6.41"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "'0.330"
$ws.Range("E12").Value = "  -3.51%  "

$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "3.043.12"
$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").Value = "58.921.05"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "'20.41"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D19").Value = "'343.72"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "'4.30"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").Value = "'10.04"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("D22").Value = "'6.40"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'66.61"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("D26").Value = "'0.404"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "'7.04"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").Value = "  -3.95%  "

$ws.Range("D31").Value = "'1.60"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").Value = "'5.88"
$ws.Range("E32").Value = "  -4.07%  "

$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").Value = "'148.99"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("D35").Value = "'3.94"
$ws.Range("E35").Value = "  -1.62%  "

$ws.Range("D36").Value = "'1.11"
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").Value = "'36.68"
$ws.Range("E37").Value = "  +1.65%  "

$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = "  -5.44%  "

$ws.Range("D40").Value = "'0.808"
$ws.Range("E40").Value = "  -6.74%  "

$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "'0.599"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").Value = "'10.75"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").Value = "'267.55"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").Value = "'0.0512"
$ws.Range("E47").Value = "  -1.91%  "

$ws.Range("D48").Value = "'18.29"
$ws.Range("E48").Value = "  -2.40%  "

$ws.Range("D49").Value = "1.958.25"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").Value = "'0.0221"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Value = "'18.08"
$ws.Range("E51").Value = "  -2.98%  "
